# Swap the contents of columns C and D (codeforiati:group-name <-> codeforiati:group-code)
# across the whole used range of the active sheet, leaving columns A and B untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colC = $ws.Range("C1:C$lastRow").Value2
$colD = $ws.Range("D1:D$lastRow").Value2

$ws.Range("C1:C$lastRow").Value2 = $colD
$ws.Range("D1:D$lastRow").Value2 = $colC
